# Auto-generated edit script: update LevePriceNQ/HQ + profit columns per scheduled market refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1978.6666
$ws.Range("I9").Value = 366.86667
$ws.Range("J9").Value = 6008.1665
$ws.Range("K9").Value = 366.86667
$ws.Range("L9").Value = 6008.1665
$ws.Range("M9").Value = -197.86667
$ws.Range("N9").Value = -6346.1665

$ws.Range("H41").Value = 1522.125
$ws.Range("I41").Value = 1574.5
$ws.Range("J41").Value = 1469.75
$ws.Range("K41").Value = 1574.5
$ws.Range("L41").Value = 1469.75
$ws.Range("M41").Value = -1134.5
$ws.Range("N41").Value = -2349.75

$ws.Range("H76").Value = 16676966
$ws.Range("J76").Value = 18975
$ws.Range("L76").Value = 18975
$ws.Range("N76").Value = -19605

$ws.Range("H79").Value = 16676966
$ws.Range("J79").Value = 18975
$ws.Range("L79").Value = 18975
$ws.Range("N79").Value = -21159

$ws.Range("H92").Value = 3463.7334
$ws.Range("I92").Value = 2639.7144
$ws.Range("K92").Value = 2639.7144
$ws.Range("M92").Value = -1391.7144

$ws.Range("H132").Value = 1407.75
$ws.Range("I132").Value = 1356.4517
$ws.Range("K132").Value = 4069.3551
$ws.Range("M132").Value = -1539.3551

$ws.Range("H138").Value = 3089.0588
$ws.Range("J138").Value = 3337.475
$ws.Range("L138").Value = 10012.425
$ws.Range("N138").Value = -20292.425

$ws.Range("H141").Value = 995
$ws.Range("I141").Value = 995
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2985
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2195
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31080.082
$ws.Range("I32").Value = 31080.082
$ws.Range("K32").Value = 31080.082
$ws.Range("M32").Value = -30793.082

$ws.Range("H45").Value = 10016
$ws.Range("J45").Value = 2766.6667
$ws.Range("L45").Value = 2766.6667
$ws.Range("N45").Value = -3520.6667

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 7028.3076
$ws.Range("I132").Value = 4436.8
$ws.Range("J132").Value = 15666.667
$ws.Range("K132").Value = 13310.4
$ws.Range("L132").Value = 47000.001
$ws.Range("M132").Value = -10780.4
$ws.Range("N132").Value = -52060.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1219.3334
$ws.Range("I20").Value = 1305.5
$ws.Range("K20").Value = 1305.5
$ws.Range("M20").Value = -1058.5

$ws.Range("H22").Value = 1482.3334
$ws.Range("I22").Value = 1630.1428
$ws.Range("J22").Value = 965
$ws.Range("K22").Value = 1630.1428
$ws.Range("L22").Value = 965
$ws.Range("M22").Value = -1457.1428
$ws.Range("N22").Value = -1311

$ws.Range("H80").Value = 80300.08
$ws.Range("J80").Value = 104191.3
$ws.Range("L80").Value = 104191.3
$ws.Range("N80").Value = -106187.3

$ws.Range("H83").Value = 80300.08
$ws.Range("J83").Value = 104191.3
$ws.Range("L83").Value = 520956.5
$ws.Range("N83").Value = -530940.5

$ws.Range("H105").Value = 20001190
$ws.Range("I105").Value = 25001072
$ws.Range("K105").Value = 25001072
$ws.Range("M105").Value = -24999325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40003656
$ws.Range("I31").Value = 55556916
$ws.Range("J31").Value = 9560
$ws.Range("K31").Value = 55556916
$ws.Range("L31").Value = 9560
$ws.Range("M31").Value = -55556621
$ws.Range("N31").Value = -10150

$ws.Range("H34").Value = 40003656
$ws.Range("I34").Value = 55556916
$ws.Range("J34").Value = 9560
$ws.Range("K34").Value = 55556916
$ws.Range("L34").Value = 9560
$ws.Range("M34").Value = -55556714
$ws.Range("N34").Value = -9964

$ws.Range("H86").Value = 6291.8335
$ws.Range("I86").Value = 5196.75
$ws.Range("K86").Value = 5196.75
$ws.Range("M86").Value = -4073.75

$ws.Range("H89").Value = 6291.8335
$ws.Range("I89").Value = 5196.75
$ws.Range("K89").Value = 25983.75
$ws.Range("M89").Value = -20367.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 100000
$ws.Range("I20").Value = 100000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 100000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -99755
$ws.Range("N20").ClearContents()

$ws.Range("H24").Value = 25490.736
$ws.Range("I24").Value = 53250
$ws.Range("J24").Value = 18088.268
$ws.Range("K24").Value = 53250
$ws.Range("L24").Value = 18088.268
$ws.Range("M24").Value = -53077
$ws.Range("N24").Value = -18434.268

$ws.Range("H49").Value = 39998.668
$ws.Range("J49").Value = 39998.668
$ws.Range("L49").Value = 39998.668
$ws.Range("N49").Value = -40366.668

$ws.Range("H55").Value = 28249
$ws.Range("I55").Value = 11515
$ws.Range("J55").Value = 44983
$ws.Range("K55").Value = 11515
$ws.Range("L55").Value = 44983
$ws.Range("M55").Value = -11188
$ws.Range("N55").Value = -45637

$ws.Range("H70").Value = 43059
$ws.Range("I70").Value = 43282
$ws.Range("K70").Value = 43282
$ws.Range("M70").Value = -43012

$ws.Range("H73").Value = 43059
$ws.Range("I73").Value = 43282
$ws.Range("K73").Value = 43282
$ws.Range("M73").Value = -42346

$ws.Range("H132").Value = 7177
$ws.Range("I132").Value = 4004.8
$ws.Range("K132").Value = 12014.4
$ws.Range("M132").Value = -9484.400000000001

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H68").Value = 14801
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251

$ws.Range("H71").Value = 14801
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256

$ws.Range("H136").Value = 5148.8965
$ws.Range("I136").Value = 3023.5789
$ws.Range("K136").Value = 9070.736699999999
$ws.Range("M136").Value = -6520.736699999999

$ws.Range("H140").Value = 51836.715
$ws.Range("J140").Value = 51836.715
$ws.Range("L140").Value = 51836.715
$ws.Range("N140").Value = -62196.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 60000
$ws.Range("I31").Value = 60000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 60000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -59652
$ws.Range("N31").ClearContents()

$ws.Range("H37").Value = 19114
$ws.Range("J37").Value = 19114
$ws.Range("L37").Value = 19114
$ws.Range("N37").Value = -19520

$ws.Range("H81").Value = 14086.182
$ws.Range("I81").Value = 2680
$ws.Range("K81").Value = 5360
$ws.Range("M81").Value = -4299

$ws.Range("H84").Value = 14086.182
$ws.Range("I84").Value = 2680
$ws.Range("K84").Value = 26800
$ws.Range("M84").Value = -21496

$ws.Range("H96").Value = 2400
$ws.Range("I96").Value = 2400
$ws.Range("K96").Value = 2400
$ws.Range("M96").Value = -1027

$ws.Range("H100").Value = 2473.9375
$ws.Range("I100").Value = 2609.3572
$ws.Range("K100").Value = 5218.7144
$ws.Range("M100").Value = -4677.7144

$ws.Range("H121").Value = 73706.164
$ws.Range("J121").Value = 73706.164
$ws.Range("L121").Value = 73706.164
$ws.Range("N121").Value = -77200.164
